$d = $word.ActiveDocument

# --- Insertion 1: new paragraph after "...within a single module." ---
# (commit: "answered question 3 and 5" -- this is the discussion of the
# proposed improvement to the step function rewards)
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("single module.", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if (-not $found1) {
    throw "Anchor text for insertion 1 not found"
}
$rng1.Collapse(0)
$ins1 = $d.Range($rng1.Start, $rng1.End)
$p1Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">A proposed improvement is implemented within </w:t></w:r><w:r><w:t>S</w:t></w:r><w:r><w:t xml:space="preserve">arsa and q-learning control implementations. The improvement </w:t></w:r><w:r><w:t>relates</w:t></w:r><w:r><w:t xml:space="preserve"> to </w:t></w:r><w:r><w:t>assigning the correct states to the step function when the AI steps on either a hole or the goal state. This was implemented to promote the improvement of the policy upon reaching goal by increasing the reward the goal state receives</w:t></w:r><w:r><w:t xml:space="preserve"> while keeping the rewards of going to holes at 0</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> Upon implementation, Sarsa and Q-Learning resulted in higher values across the board, by a factor of 10. Removing the implementation returns the expected value with both implementations returning optimal policies.</w:t></w:r></w:p>'
$ins1.InsertXML($p1Xml)

# --- Insertion 2: four new paragraphs (answers to questions 3 and 5) after the
# "...to return an optimal policy." paragraph, before the trailing blank paragraphs ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("it did require more iterations to return an optimal policy.", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if (-not $found2) {
    throw "Anchor text for insertion 2 not found"
}
$rng2.Collapse(0)
$ins2 = $d.Range($rng2.Start, $rng2.End)
$block2Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>3.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">For the purpose of this analysis, the policy move of death states </w:t></w:r><w:r><w:t>and goal state are not relevant.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>In this case, an optimal policy is given when all non-states reach the goal state.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Sarsa control achieved an optimal policy after approximately 900 states while </w:t></w:r><w:r><w:t>Q-Learning achieved an optimal policy after</w:t></w:r><w:r><w:t xml:space="preserve"> around 700 states.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>However, these evaluations vary wildly between runs</w:t></w:r><w:r><w:t xml:space="preserve"> with the common point being the lack of a policy for the corner state 3 with rewards given. </w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>Evaluation of Sarsa</w:t></w:r><w:r><w:t xml:space="preserve"> control </w:t></w:r><w:r><w:t>(w/ proposed improvements) returned an optimal policy after approximately 300 episodes.</w:t></w:r><w:r><w:t xml:space="preserve"> While Q-learning control (w/ proposed improvements) returned an optimal policy for all states after approximately 600 episodes. On average, Sarsa</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">performs better than Q-learning, with the proposed improvements decreasing the episodes required for an optimal policy by 250 episodes on average. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">5. After numerous hyperparameter changes, it was </w:t></w:r><w:r><w:t>not</w:t></w:r><w:r><w:t xml:space="preserve"> possible to return an optimal policy for the big frozen lake</w:t></w:r><w:r><w:t xml:space="preserve"> for every stat</w:t></w:r><w:r><w:t>e</w:t></w:r><w:r><w:t xml:space="preserve"> for Sarsa</w:t></w:r><w:r><w:t xml:space="preserve">. This was regardless of </w:t></w:r><w:r><w:t xml:space="preserve">maximum </w:t></w:r><w:r><w:t>Eps</w:t></w:r><w:r><w:t xml:space="preserve">ilon and </w:t></w:r><w:r><w:t>T</w:t></w:r><w:r><w:t xml:space="preserve">heta to promote exploration, with </w:t></w:r><w:r><w:t>an exceedingly large number of episodes</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>It was</w:t></w:r><w:r><w:t xml:space="preserve"> able</w:t></w:r><w:r><w:t xml:space="preserve"> to return a policy</w:t></w:r><w:r><w:t xml:space="preserve"> that defines a </w:t></w:r><w:r><w:t xml:space="preserve">safe path towards the goal. However, the model was too </w:t></w:r><w:r><w:t>conservative</w:t></w:r><w:r><w:t xml:space="preserve"> and was not able to explore much of the bottom left of the large lake as opposed to the top right of the lake</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p><w:p><w:r><w:t>This is because Sarsa is an on-policy algorithm which will follow the policy to compute the next state. This is not the case with Q learning which was able to return an optimal policy with much less states</w:t></w:r><w:r><w:t xml:space="preserve"> of around 20,000 states. This is because QL takes the maximum reward of the new state and ignores the current policy, enabling a more volatile but with more exploration.</w:t></w:r></w:p>'
$ins2.InsertXML($block2Xml)

# The "5. After numerous..." and "This is because Sarsa..." paragraphs use an
# explicit zero indent (left=0, firstLine=0) in the target document, so set
# that directly on the paragraph format (InsertXML drops an explicit
# firstLine="0" since it is already the implicit default).
$rngA = $d.Content
$foundA = $rngA.Find.Execute("5. After numerous hyperparameter changes", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if (-not $foundA) {
    throw "Paragraph 5 not found for indent fixup"
}
$paraA = $rngA.Paragraphs.Item(1)
$paraA.LeftIndent = 0
$paraA.FirstLineIndent = 0

$rngB = $d.Content
$foundB = $rngB.Find.Execute("This is because Sarsa is an on-policy algorithm", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if (-not $foundB) {
    throw "Trailing paragraph not found for indent fixup"
}
$paraB = $rngB.Paragraphs.Item(1)
$paraB.LeftIndent = 0
$paraB.FirstLineIndent = 0
